$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Timestamp update (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 19 de Junio de 2020 a las 17:01"

# --- Updated stats for existing countries (rows unchanged) ---
# Row 4 - Estados Unidos
$ws.Range("B4").Value = 2265449
$ws.Range("C4").Value = 1798
$ws.Range("D4").Value = 931355
$ws.Range("E4").Value = 1213368
$ws.Range("G4").Value = 38
$ws.Range("H4").Value = 120726

# Row 7 - India
$ws.Range("B7").Value = 384953
$ws.Range("C7").Value = 3862
$ws.Range("E7").Value = 166442
$ws.Range("G7").Value = 63
$ws.Range("H7").Value = 12667

# Row 8 - Reino Unido
$ws.Range("B8").Value = 301815
$ws.Range("C8").Value = 1346
$ws.Range("G8").Value = 173
$ws.Range("H8").Value = 42461

# Row 14 - Alemania
$ws.Range("B14").Value = 190290
$ws.Range("C14").Value = 164
$ws.Range("E14").Value = 7238

# Row 34 - Paises Bajos
$ws.Range("D34").Value = 33459
$ws.Range("E34").Value = 8130

# Row 37 - Emiratos Arabes Unidos
$ws.Range("D37").Value = 11851
$ws.Range("E37").Value = 24705
$ws.Range("G37").Value = 6
$ws.Range("H37").Value = 954

# Row 46 - (#50)
$ws.Range("B46").Value = 25068
$ws.Range("C46").Value = 423
$ws.Range("D46").Value = 14605
$ws.Range("E46").Value = 9816
$ws.Range("G46").Value = 12
$ws.Range("H46").Value = 647

# Row 79 - (#83)
$ws.Range("B79").Value = 5338
$ws.Range("C79").Value = 59
$ws.Range("D79").Value = 3830
$ws.Range("E79").Value = 1457

# Row 82 - (#86)
$ws.Range("B82").Value = 4820
$ws.Range("C82").Value = 156
$ws.Range("D82").Value = 1863
$ws.Range("E82").Value = 2735
$ws.Range("G82").Value = 6
$ws.Range("H82").Value = 222

# Row 100 - (#104)
$ws.Range("B100").Value = 2305
$ws.Range("C100").Value = 10
$ws.Range("D100").Value = 2037
$ws.Range("E100").Value = 183

# Row 130 - (#134)
$ws.Range("B130").Value = 900
$ws.Range("C130").Value = 1
$ws.Range("E130").Value = 37

# --- Reordered country rows (names swapped with new stats) ---
# Rows 182/183: Liechtenstein <-> Botsuana
$ws.Range("A182").Value = "Botsuana"
$ws.Range("B182").Value = 89
$ws.Range("C182").Value = 10
$ws.Range("D182").Value = 25
$ws.Range("E182").Value = 63

$ws.Range("A183").Value = "Liechtenstein"
$ws.Range("B183").Value = 82
$ws.Range("D183").Value = 55
$ws.Range("E183").Value = 26

# Rows 202/203: Fiyi <-> Dominica
$ws.Range("A202").Value = "Dominica"
$ws.Range("A203").Value = "Fiyi"

# Rows 206/207: Groenlandia <-> Islas Malvinas
$ws.Range("A206").Value = "Islas Malvinas"
$ws.Range("A207").Value = "Groenlandia"

# Rows 210/211: Montserrat <-> Seychelles
$ws.Range("A210").Value = "Seychelles"
$ws.Range("D210").Value = 11
$ws.Range("H210").Value = 0

$ws.Range("A211").Value = "Montserrat"
$ws.Range("D211").Value = 10
$ws.Range("H211").Value = 1

# Rows 213/214: Islas Virgenes Britanicas <-> Papua Nueva Guinea
$ws.Range("A213").Value = "Papua Nueva Guinea"
$ws.Range("D213").Value = 8
$ws.Range("H213").Value = 0

$ws.Range("A214").Value = "Islas Virgenes Britanicas"
$ws.Range("D214").Value = 7
$ws.Range("H214").Value = 1
